$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.17413792331684164,
    -0.0059999999546711535,
    -0.0039999999594417801,
    -0.0079999999258983934,
    -0.0029999999577299263,
    0.052352782022124345,
    -0.0099999998952768898,
    -0.052095558227846617,
    -0.0019999999535080804,
    -0.0019999999569364491,
    0.047969206781237617,
    -0.0034999999444065288,
    -0.0034999999402502979,
    -0.0079999999064597205,
    -0.00099999995684285636,
    -0.0019999999489197506,
    -0.0019999999480448949,
    -0.0039999999331978842,
    -0.0039999999689812604,
    -0.0039999999665436548,
    -0.0039999999661759489,
    -0.0039999999659006136,
    -0.0049999999495131675,
    -0.019999999836064042,
    -0.019999999833713034,
    -0.0024999999488084512,
    -0.0024999999455777022,
    -0.0019999999335871266,
    -0.0069999998868199853,
    -0.059999999495193634,
    -0.0069999998776797412,
    -0.0099999998553048641,
    -0.0039999998988697882
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
